$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H16").Value = 353
$ws_ALC.Range("I16").Value = 179.5
$ws_ALC.Range("J16").Value = 700
$ws_ALC.Range("K16").Value = 179.5
$ws_ALC.Range("L16").Value = 700
$ws_ALC.Range("M16").Value = 50.5

$ws_ALC.Range("H33").Value = 1030.6364
$ws_ALC.Range("I33").Value = 292.125
$ws_ALC.Range("J33").Value = 3000
$ws_ALC.Range("K33").Value = 292.125
$ws_ALC.Range("L33").Value = 3000
$ws_ALC.Range("M33").Value = -63.125
$ws_ALC.Range("N33").Value = -3458

$ws_ALC.Range("H38").Value = 11924.4
$ws_ALC.Range("I38").Value = 12360.444
$ws_ALC.Range("J38").Value = 8000
$ws_ALC.Range("K38").Value = 37081.33199999999
$ws_ALC.Range("L38").Value = 24000
$ws_ALC.Range("M38").Value = -36709.33199999999

$ws_ALC.Range("H40").Value = 3201
$ws_ALC.Range("I40").Value = 2400
$ws_ALC.Range("J40").Value = 4002
$ws_ALC.Range("K40").Value = 2400
$ws_ALC.Range("L40").Value = 4002
$ws_ALC.Range("M40").Value = -2225
$ws_ALC.Range("N40").Value = -4352

$ws_ALC.Range("H98").Value = 2375
$ws_ALC.Range("I98").Value = 2375
$ws_ALC.Range("J98").Value = 0
$ws_ALC.Range("K98").Value = 2375
$ws_ALC.Range("L98").Value = 0
$ws_ALC.Range("M98").Value = -877
$ws_ALC.Range("N98").ClearContents()

$ws_ALC.Range("H100").Value = 1049.25
$ws_ALC.Range("I100").Value = 999
$ws_ALC.Range("J100").Value = 1200
$ws_ALC.Range("K100").Value = 999
$ws_ALC.Range("L100").Value = 1200
$ws_ALC.Range("M100").Value = -458

$ws_ALC.Range("H112").Value = 1358.5
$ws_ALC.Range("I112").Value = 0
$ws_ALC.Range("J112").Value = 1358.5
$ws_ALC.Range("K112").Value = 0
$ws_ALC.Range("L112").Value = 4075.5
$ws_ALC.Range("N112").Value = -6291.5

$ws_ALC.Range("H119").Value = 0
$ws_ALC.Range("I119").Value = 0
$ws_ALC.Range("J119").Value = 0
$ws_ALC.Range("K119").Value = 0
$ws_ALC.Range("L119").Value = 0
$ws_ALC.Range("N119").ClearContents()

$ws_ALC.Range("H122").Value = 2375
$ws_ALC.Range("I122").Value = 2375
$ws_ALC.Range("J122").Value = 0
$ws_ALC.Range("K122").Value = 7125
$ws_ALC.Range("L122").Value = 0
$ws_ALC.Range("M122").Value = -4675
$ws_ALC.Range("N122").ClearContents()

$ws_ARM.Range("H31").Value = 7333.3335
$ws_ARM.Range("I31").Value = 7333.3335
$ws_ARM.Range("J31").Value = 0
$ws_ARM.Range("K31").Value = 7333.3335
$ws_ARM.Range("L31").Value = 0
$ws_ARM.Range("M31").Value = -7039.3335

$ws_ARM.Range("H45").Value = 225
$ws_ARM.Range("I45").Value = 225
$ws_ARM.Range("J45").Value = 0
$ws_ARM.Range("K45").Value = 225
$ws_ARM.Range("L45").Value = 0
$ws_ARM.Range("M45").Value = 152

$ws_ARM.Range("H61").Value = 1000.4
$ws_ARM.Range("I61").Value = 1000.4
$ws_ARM.Range("J61").Value = 0
$ws_ARM.Range("K61").Value = 1000.4
$ws_ARM.Range("L61").Value = 0
$ws_ARM.Range("M61").Value = -788.4

$ws_ARM.Range("H110").Value = 1293.125
$ws_ARM.Range("I110").Value = 1356.5714
$ws_ARM.Range("J110").Value = 849
$ws_ARM.Range("K110").Value = 1356.5714
$ws_ARM.Range("L110").Value = 849
$ws_ARM.Range("M110").Value = 688.4286

$ws_ARM.Range("H136").Value = 1000.4
$ws_ARM.Range("I136").Value = 1000.4
$ws_ARM.Range("J136").Value = 0
$ws_ARM.Range("K136").Value = 3001.2
$ws_ARM.Range("L136").Value = 0
$ws_ARM.Range("M136").Value = -451.1999999999998

$ws_BSM.Range("H25").Value = 514
$ws_BSM.Range("I25").Value = 514
$ws_BSM.Range("J25").Value = 0
$ws_BSM.Range("K25").Value = 514
$ws_BSM.Range("L25").Value = 0
$ws_BSM.Range("M25").Value = -279

$ws_BSM.Range("H94").Value = 1050
$ws_BSM.Range("I94").Value = 1050
$ws_BSM.Range("J94").Value = 0
$ws_BSM.Range("K94").Value = 1050
$ws_BSM.Range("L94").Value = 0
$ws_BSM.Range("M94").Value = -599

$ws_BSM.Range("H96").Value = 10498.333
$ws_BSM.Range("I96").Value = 10498.333
$ws_BSM.Range("J96").Value = 0
$ws_BSM.Range("K96").Value = 10498.333
$ws_BSM.Range("L96").Value = 0
$ws_BSM.Range("M96").Value = -7752.333000000001

$ws_BSM.Range("H97").Value = 22558.5
$ws_BSM.Range("I97").Value = 26462.75
$ws_BSM.Range("J97").Value = 14750
$ws_BSM.Range("K97").Value = 26462.75
$ws_BSM.Range("L97").Value = 14750
$ws_BSM.Range("M97").Value = -25471.75
$ws_BSM.Range("N97").Value = -16732

$ws_BSM.Range("H99").Value = 3125
$ws_BSM.Range("I99").Value = 3125
$ws_BSM.Range("J99").Value = 0
$ws_BSM.Range("K99").Value = 3125
$ws_BSM.Range("L99").Value = 0
$ws_BSM.Range("M99").Value = -1627

$ws_BSM.Range("H102").Value = 7250
$ws_BSM.Range("I102").Value = 7250
$ws_BSM.Range("J102").Value = 0
$ws_BSM.Range("K102").Value = 7250
$ws_BSM.Range("L102").Value = 0
$ws_BSM.Range("M102").Value = -4005

$ws_BSM.Range("H105").Value = 1900
$ws_BSM.Range("I105").Value = 0
$ws_BSM.Range("J105").Value = 1900
$ws_BSM.Range("K105").Value = 0
$ws_BSM.Range("L105").Value = 1900
$ws_BSM.Range("M105").ClearContents()
$ws_BSM.Range("N105").Value = -5394

$ws_BSM.Range("H134").Value = 1000
$ws_BSM.Range("I134").Value = 1000
$ws_BSM.Range("J134").Value = 0
$ws_BSM.Range("K134").Value = 3000
$ws_BSM.Range("L134").Value = 0
$ws_BSM.Range("M134").Value = -465

$ws_CRP.Range("H2").Value = 17724.75
$ws_CRP.Range("I2").Value = 70000
$ws_CRP.Range("J2").Value = 299.66666
$ws_CRP.Range("K2").Value = 70000
$ws_CRP.Range("L2").Value = 299.66666
$ws_CRP.Range("M2").Value = -69887
$ws_CRP.Range("N2").Value = -525.66666

$ws_CRP.Range("H23").Value = 1755

$ws_CRP.Range("H27").Value = 1755

$ws_CRP.Range("H63").Value = 20000
$ws_CRP.Range("I63").Value = 0
$ws_CRP.Range("J63").Value = 20000
$ws_CRP.Range("K63").Value = 0
$ws_CRP.Range("L63").Value = 20000
$ws_CRP.Range("N63").Value = -21372

$ws_CRP.Range("H66").Value = 20000
$ws_CRP.Range("I66").Value = 0
$ws_CRP.Range("J66").Value = 20000
$ws_CRP.Range("K66").Value = 0
$ws_CRP.Range("L66").Value = 60000
$ws_CRP.Range("N66").Value = -66864

$ws_CRP.Range("H93").Value = 15000
$ws_CRP.Range("I93").Value = 15000
$ws_CRP.Range("J93").Value = 0
$ws_CRP.Range("K93").Value = 15000
$ws_CRP.Range("L93").Value = 0
$ws_CRP.Range("M93").Value = -13128

$ws_CRP.Range("H105").Value = 1500
$ws_CRP.Range("I105").Value = 1500
$ws_CRP.Range("J105").Value = 0
$ws_CRP.Range("K105").Value = 1500
$ws_CRP.Range("L105").Value = 0
$ws_CRP.Range("M105").Value = 247

$ws_CUL.Range("H81").Value = 0
$ws_CUL.Range("I81").Value = 0
$ws_CUL.Range("J81").Value = 0
$ws_CUL.Range("K81").Value = 0
$ws_CUL.Range("L81").Value = 0
$ws_CUL.Range("M81").ClearContents()

$ws_CUL.Range("H84").Value = 0
$ws_CUL.Range("I84").Value = 0
$ws_CUL.Range("J84").Value = 0
$ws_CUL.Range("K84").Value = 0
$ws_CUL.Range("L84").Value = 0
$ws_CUL.Range("M84").ClearContents()

$ws_CUL.Range("H97").Value = 787.5
$ws_CUL.Range("I97").Value = 1000
$ws_CUL.Range("J97").Value = 716.6667
$ws_CUL.Range("K97").Value = 3000
$ws_CUL.Range("L97").Value = 2150.0001
$ws_CUL.Range("M97").Value = -2504
$ws_CUL.Range("N97").Value = -3142.0001

$ws_CUL.Range("H109").Value = 1999
$ws_CUL.Range("I109").Value = 1999
$ws_CUL.Range("J109").Value = 0
$ws_CUL.Range("K109").Value = 5997
$ws_CUL.Range("L109").Value = 0
$ws_CUL.Range("M109").Value = -4957

$ws_CUL.Range("H137").Value = 7016.5
$ws_CUL.Range("I137").Value = 0
$ws_CUL.Range("J137").Value = 7016.5
$ws_CUL.Range("K137").Value = 0
$ws_CUL.Range("L137").Value = 21049.5
$ws_CUL.Range("N137").Value = -31249.5

$ws_GSM.Range("H10").Value = 2249.75
$ws_GSM.Range("I10").Value = 1999.6666
$ws_GSM.Range("J10").Value = 3000
$ws_GSM.Range("K10").Value = 1999.6666
$ws_GSM.Range("L10").Value = 3000
$ws_GSM.Range("M10").Value = -1830.6666
$ws_GSM.Range("N10").Value = -3338

$ws_GSM.Range("H17").Value = 5000
$ws_GSM.Range("I17").Value = 5000
$ws_GSM.Range("J17").Value = 0
$ws_GSM.Range("K17").Value = 5000
$ws_GSM.Range("L17").Value = 0
$ws_GSM.Range("M17").Value = -4832
$ws_GSM.Range("N17").ClearContents()

$ws_GSM.Range("H19").Value = 1100
$ws_GSM.Range("I19").Value = 1200
$ws_GSM.Range("J19").Value = 1000
$ws_GSM.Range("K19").Value = 1200
$ws_GSM.Range("L19").Value = 1000
$ws_GSM.Range("M19").Value = -912
$ws_GSM.Range("N19").Value = -1576

$ws_GSM.Range("H102").Value = 1350
$ws_GSM.Range("I102").Value = 1350
$ws_GSM.Range("J102").Value = 0
$ws_GSM.Range("K102").Value = 1350
$ws_GSM.Range("L102").Value = 0
$ws_GSM.Range("M102").Value = 272

$ws_GSM.Range("H113").Value = 4595.5557
$ws_GSM.Range("I113").Value = 4595.5557
$ws_GSM.Range("J113").Value = 0
$ws_GSM.Range("K113").Value = 4595.5557
$ws_GSM.Range("L113").Value = 0
$ws_GSM.Range("M113").Value = -2425.5557

$ws_GSM.Range("H122").Value = 4399.2
$ws_GSM.Range("I122").Value = 4998.6665
$ws_GSM.Range("J122").Value = 3500
$ws_GSM.Range("K122").Value = 14995.9995
$ws_GSM.Range("L122").Value = 10500
$ws_GSM.Range("M122").Value = -12545.9995
$ws_GSM.Range("N122").Value = -15400

$ws_LTW.Range("H17").Value = 3342.3333
$ws_LTW.Range("I17").Value = 18
$ws_LTW.Range("J17").Value = 5004.5
$ws_LTW.Range("K17").Value = 18
$ws_LTW.Range("L17").Value = 5004.5
$ws_LTW.Range("M17").Value = 152
$ws_LTW.Range("N17").Value = -5344.5

$ws_LTW.Range("H31").Value = 25374.75
$ws_LTW.Range("I31").Value = 2750
$ws_LTW.Range("J31").Value = 47999.5
$ws_LTW.Range("K31").Value = 2750
$ws_LTW.Range("L31").Value = 47999.5
$ws_LTW.Range("M31").Value = -2502
$ws_LTW.Range("N31").Value = -48495.5

$ws_LTW.Range("H40").Value = 7174.2
$ws_LTW.Range("I40").Value = 7174.2
$ws_LTW.Range("J40").Value = 0
$ws_LTW.Range("K40").Value = 7174.2
$ws_LTW.Range("L40").Value = 0
$ws_LTW.Range("M40").Value = -7038.2
$ws_LTW.Range("N40").ClearContents()

$ws_LTW.Range("H61").Value = 1000
$ws_LTW.Range("I61").Value = 1000
$ws_LTW.Range("J61").Value = 0
$ws_LTW.Range("K61").Value = 1000
$ws_LTW.Range("L61").Value = 0
$ws_LTW.Range("M61").Value = -798

$ws_LTW.Range("H68").Value = 2000
$ws_LTW.Range("I68").Value = 2000
$ws_LTW.Range("J68").Value = 0
$ws_LTW.Range("K68").Value = 2000
$ws_LTW.Range("L68").Value = 0
$ws_LTW.Range("M68").Value = -1251

$ws_LTW.Range("H71").Value = 2000
$ws_LTW.Range("I71").Value = 2000
$ws_LTW.Range("J71").Value = 0
$ws_LTW.Range("K71").Value = 10000
$ws_LTW.Range("L71").Value = 0
$ws_LTW.Range("M71").Value = -6256

$ws_LTW.Range("H113").Value = 1000
$ws_LTW.Range("I113").Value = 1000
$ws_LTW.Range("J113").Value = 0
$ws_LTW.Range("K113").Value = 1000
$ws_LTW.Range("L113").Value = 0
$ws_LTW.Range("M113").Value = 1170

$ws_LTW.Range("H122").Value = 4947.3335
$ws_LTW.Range("I122").Value = 4936.8
$ws_LTW.Range("J122").Value = 5000
$ws_LTW.Range("K122").Value = 14810.4
$ws_LTW.Range("L122").Value = 15000
$ws_LTW.Range("M122").Value = -12360.4

$ws_WVR.Range("H3").Value = 3626500.8
$ws_WVR.Range("I3").Value = 3626500.8
$ws_WVR.Range("J3").Value = 0
$ws_WVR.Range("K3").Value = 3626500.8
$ws_WVR.Range("L3").Value = 0
$ws_WVR.Range("M3").Value = -3626386.8

$ws_WVR.Range("H6").Value = 0
$ws_WVR.Range("I6").Value = 0
$ws_WVR.Range("J6").Value = 0
$ws_WVR.Range("K6").Value = 0
$ws_WVR.Range("L6").Value = 0
$ws_WVR.Range("M6").ClearContents()

$ws_WVR.Range("H11").Value = 10000000
$ws_WVR.Range("I11").Value = 10000000
$ws_WVR.Range("J11").Value = 0
$ws_WVR.Range("K11").Value = 10000000
$ws_WVR.Range("L11").Value = 0
$ws_WVR.Range("M11").Value = -9999858

$ws_WVR.Range("H113").Value = 2266.6667
$ws_WVR.Range("I113").Value = 900
$ws_WVR.Range("J113").Value = 2950
$ws_WVR.Range("K113").Value = 2700
$ws_WVR.Range("L113").Value = 8850
$ws_WVR.Range("M113").Value = -530
$ws_WVR.Range("N113").Value = -13190

$ws_WVR.Range("H122").Value = 2584.8572
$ws_WVR.Range("I122").Value = 2584.8572
$ws_WVR.Range("J122").Value = 0
$ws_WVR.Range("K122").Value = 7754.571599999999
$ws_WVR.Range("L122").Value = 0
$ws_WVR.Range("M122").Value = -5304.571599999999

$ws_WVR.Range("H132").Value = 6779.6665
$ws_WVR.Range("I132").Value = 7135.8
$ws_WVR.Range("J132").Value = 4999
$ws_WVR.Range("K132").Value = 21407.4
$ws_WVR.Range("L132").Value = 14997
$ws_WVR.Range("M132").Value = -20058.5
$ws_WVR.Range("N132").Value = -20057
